# Append 4 new incident-log rows (132-135) to the bottom of the "tickets"
# sheet, matching the newly logged tickets from 2024-05-21.
#
# Column A holds plain-text dates (e.g. "2024-05-21"). Assigning a bare
# date-shaped string through .Value makes Excel auto-convert it into a
# real date serial, which is not what the source data wants (every
# existing cell in the sheet is stored as literal text). Prefixing the
# value with a leading apostrophe forces Excel to keep it as literal text
# (the standard "format as text" entry trick), while every other column
# here (times, durations, dashes, free text) round-trips as text already
# without any extra help.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(132, "2024-05-21", "11:50:30", "Fallo tornillo",              "-", "-", "-", "-", "11:50:32", "0:00:02"),
    @(133, "2024-05-21", "11:50:35", "Palet atascado en la curva",  "-", "-", "-", "-", "11:50:42", "0:00:07"),
    @(134, "2024-05-21", "11:50:39", "Fallo tornillo",              "-", "-", "-", "-", "11:50:40", "0:00:01"),
    @(135, "2024-05-21", "11:51:29", "Fallo tornillo",              "-", "-", "-", "-", "11:51:31", "0:00:02")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = "'" + $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Range("G$r").Value = $row[7]
    $ws.Range("H$r").Value = $row[8]
    $ws.Range("I$r").Value = $row[9]
}
